# "Generate Report for Handback"
#
# The file "abb6b744-1871-4545-a88f-6fac0f1d6f99.md" has now been handed
# back successfully (previously it was "Ready for handoff" / showed a
# stale-handback error). Update the localization status report to
# reflect the new handback state across all three sheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet ------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Handed back: in sync with en-US"
$zhcn.Range("K3").Value = "2016-08-26 12:49:07"
$zhcn.Range("P3").Value = ""

# --- de-de sheet --------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Handed back: in sync with en-US"
$dede.Range("K3").Value = "2016-08-26 12:49:16"
$dede.Range("P3").Value = ""
